$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.15   # Current Capital
$summary.Range("B4").Value = 0.95      # Total P&L $
$summary.Range("B5").Value = 0.21      # Total P&L %
$summary.Range("B6").Value = 90        # Total Trades
$summary.Range("B7").Value = 44        # Winning Trades
$summary.Range("B9").Value = 48.89     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: Strategy Status (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.15     # Capital
$status.Range("D5").Value = 57         # Trades
$status.Range("E5").Value = 0.84       # P&L $
$status.Range("F5").Value = 1.15       # P&L %
$status.Range("G5").Value = 52.63      # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: All Trades - Trade #90 (row 91) closes out
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G91").Value = 0.85
$allTrades.Range("H91").Value = "CLOSED"
$allTrades.Range("I91").Value = 10.3896
$allTrades.Range("J91").Value = 0.08
$allTrades.Range("K91").Value = 101.15
$allTrades.Range("L91").Value = "early_exit"
$allTrades.Range("M91").Value = 0.13

# ---------------------------------------------------------------------------
# Sheet: All Trades - new Trade #123 (row 124)
# ---------------------------------------------------------------------------
$allTrades.Range("A124").Value = 123
$allTrades.Range("B124").Value = "'2026-02-17"
$allTrades.Range("C124").Value = "21:11:05"
$allTrades.Range("D124").Value = "MarketMaking"
$allTrades.Range("E124").Value = "UP"
$allTrades.Range("F124").Value = 0.77
$allTrades.Range("G124").Value = ""
$allTrades.Range("H124").Value = "OPEN"
$allTrades.Range("I124").Value = 0
$allTrades.Range("J124").Value = 0
$allTrades.Range("K124").Value = 101.0746450978375
$allTrades.Range("L124").Value = ""
$allTrades.Range("M124").Value = 0
$allTrades.Range("N124").Value = 0
$allTrades.Range("O124").Value = 0
$allTrades.Range("P124").Value = 0.6
$allTrades.Range("Q124").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Sheet: MarketMaking - Trade #90 (row 58) closes out
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G58").Value = 0.85
$mm.Range("H58").Value = "CLOSED"
$mm.Range("I58").Value = 10.3896
$mm.Range("J58").Value = 0.08
$mm.Range("K58").Value = 101.15
$mm.Range("P58").Value = "early_exit"
$mm.Range("Q58").Value = 0.13

# ---------------------------------------------------------------------------
# Sheet: MarketMaking - new Trade #123 (row 91)
# ---------------------------------------------------------------------------
$mm.Range("A91").Value = 123
$mm.Range("B91").Value = "'2026-02-17"
$mm.Range("C91").Value = "21:11:05"
$mm.Range("D91").Value = "MarketMaking"
$mm.Range("E91").Value = "UP"
$mm.Range("F91").Value = 0.77
$mm.Range("G91").Value = ""
$mm.Range("H91").Value = "OPEN"
$mm.Range("I91").Value = 0
$mm.Range("J91").Value = 0
$mm.Range("K91").Value = 101.0746450978375
$mm.Range("L91").Value = 0
$mm.Range("M91").Value = 0
$mm.Range("N91").Value = 0.6
$mm.Range("O91").Value = "Normal spread capture: 19600 bps"
$mm.Range("P91").Value = ""
$mm.Range("Q91").Value = 0
